# Start FFT overlap implementation
# Update the "prescaler" requirement for TIM4 (column C) from 100 to 1000,
# which in turn recalculates dependent formula cells (C5, C7, C12),
# and move the active cell selection to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 1000

$ws.Range("C4").Select()

$wb.Save()
